$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("Create a country", "PASSED", "chrome"),
    @("Create a country", "PASSED", "chrome"),
    @("Create a country 2", "PASSED", "chrome"),
    @("Create a Citizenship", "PASSED", "chrome"),
    @("Login with valid username and password", "PASSED", "chrome"),
    @("Create a country", "PASSED", "chrome"),
    @("Create a country", "PASSED", "chrome"),
    @("Create a country 2", "PASSED", "chrome"),
    @("Create a citizenship", "PASSED", "chrome")
)

$startRow = 22
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
}
